$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores Price/Volume/Hora as literal text (inline strings), not
# numbers -- e.g. "319.72", "-3.62%", "12" are text labels scraped from a
# web page. Force each touched cell to Text format before writing so Excel
# keeps the new content as a literal string instead of silently coercing it
# into a number/percentage (which would change both the stored type and the
# value itself, e.g. "-3.62%" -> -0.0362).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "319.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.62%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "12"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "42.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.52%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "12"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-5.37%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "12"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08155"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.50%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "12"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.324"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.20%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "12"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.797"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-14.31%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "12"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9357"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.45%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "12"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1107"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.50%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "12"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1853"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.15%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "12"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09376"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.52%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "12"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04616"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.55%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "12"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.425"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-21.40%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "12"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "12"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.21%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "12"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005703"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.28%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "12"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.357"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.01%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "12"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.544"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.63%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "12"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3348"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.65%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "12"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1381"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.02%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "12"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2523"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.87%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "12"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04155"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.20%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "12"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001247"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.64%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "12"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.53%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "12"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001202"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-7.68%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "12"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002982"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.38%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "12"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "12"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "12"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "12"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "12"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "12"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "12"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "12"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "12"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "12"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "12"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "12"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02705"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-0.63%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "12"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05541"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.94%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "12"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008044"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.93%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "12"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1397"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.70%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "12"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006548"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-12.45%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "12"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002082"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.01%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "12"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007612"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.93%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "12"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3174"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.45%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "12"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006934"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.90%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "12"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "12"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003333"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "9.10%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "12"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.16%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "12"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "12"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "12"
